$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# These columns store numeric-looking figures as plain text (shared strings),
# so force text format before assigning to avoid Excel auto-converting them
# to numbers.
$cells = @("B11", "C11", "C12", "D12")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Enterprises density (per 1000 people) row: Micro / SMEs values
$ws.Range("B11").Value = "11.85"
$ws.Range("C11").Value = "10.45"

# Enterprises (% of total) row: SMEs / MSMEs values
$ws.Range("C12").Value = "46.13"
$ws.Range("D12").Value = "98.43"
